# Update marksheet correct/total marks figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking row, Right column: 3 -> 5
$ws.Range("B11").Value = 5

# B12: Total row, Right column: 81 -> 135
$ws.Range("B12").Value = 135

# E12: Total row, Max column (text "corr/total"): "80/84" -> "135/140"
$ws.Range("E12").Value = "135/140"
